# "Generate Report for Handoff"
#
# The localization-status report is regenerated: the file
# "3976662b-7e18-4ef4-9397-f67d8c4c4bc5" is handed off again (status goes
# from "Handed back: in sync with en-US" to "Ready for handoff", and a new
# Latest Handoff Datetime is recorded), while
# "efbf9219-8bb8-4182-8079-84878fca9520" keeps its previous
# "Handed back: in sync with en-US" status. During report (re)generation
# the two files swap table rows (efbf9219 now sorts into row 2, 3976662b
# into row 3) on every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "efbf9219-8bb8-4182-8079-84878fca9520.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"

$ws.Range("A3").Value = "3976662b-7e18-4ef4-9397-f67d8c4c4bc5.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"

$ws.Range("A4").Value = ".localization-config"
$ws.Range("B4").Value = "Not to be localized"
$ws.Range("C4").Value = "Not to be localized"

# Rebuild the hyperlinks, keeping the same Address per r:id slot as before
# (the relationship targets are untouched by this edit) but pointing the
# slots at the rows they now belong to.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/74cfad7bf15d4feed10d948ffefb164bb9c05f74/e2e/3976662b-7e18-4ef4-9397-f67d8c4c4bc5.md", "", "", "efbf9219-8bb8-4182-8079-84878fca9520.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/74cfad7bf15d4feed10d948ffefb164bb9c05f74/e2e/efbf9219-8bb8-4182-8079-84878fca9520.md", "", "", "3976662b-7e18-4ef4-9397-f67d8c4c4bc5.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/74cfad7bf15d4feed10d948ffefb164bb9c05f74/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "efbf9219-8bb8-4182-8079-84878fca9520.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "efbf9219-8bb8-4182-8079-84878fca9520.d345dcc8bf0012c19232a2d193a32928fc9c98ef.zh-cn.xlf"
$ws.Range("D2").Value = "2016-03-08 20:53:21"
$ws.Range("E2").Value = "efbf9219-8bb8-4182-8079-84878fca9520.md"
$ws.Range("F2").Value = "efbf9219-8bb8-4182-8079-84878fca9520.d345dcc8bf0012c19232a2d193a32928fc9c98ef.zh-cn.xlf"
$ws.Range("G2").Value = "2016-03-08 20:54:29"
$ws.Range("H2").Value = "Include"

$ws.Range("A3").Value = "3976662b-7e18-4ef4-9397-f67d8c4c4bc5.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "3976662b-7e18-4ef4-9397-f67d8c4c4bc5.7fed2e4b03b1d9d3d15c2b653214b755747046ba.zh-cn.xlf"
$ws.Range("D3").Value = "2016-03-08 20:55:43"
$ws.Range("E3").Value = "3976662b-7e18-4ef4-9397-f67d8c4c4bc5.md"
$ws.Range("F3").Value = "3976662b-7e18-4ef4-9397-f67d8c4c4bc5.7fed2e4b03b1d9d3d15c2b653214b755747046ba.zh-cn.xlf"
$ws.Range("G3").Value = "2016-03-08 20:54:29"
$ws.Range("H3").Value = "Include"

$ws.Range("A4").Value = ".localization-config"
$ws.Range("B4").Value = "Not to be localized"
$ws.Range("D4").Value = "0001-01-01 00:00:00"
$ws.Range("G4").Value = "0001-01-01 00:00:00"
$ws.Range("H4").Value = "Ignored"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/74cfad7bf15d4feed10d948ffefb164bb9c05f74/e2e/3976662b-7e18-4ef4-9397-f67d8c4c4bc5.md", "", "", "efbf9219-8bb8-4182-8079-84878fca9520.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/988e7166473fb56d35401be470c3d6e8703d49a9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3976662b-7e18-4ef4-9397-f67d8c4c4bc5.7fed2e4b03b1d9d3d15c2b653214b755747046ba.zh-cn.xlf", "", "", "efbf9219-8bb8-4182-8079-84878fca9520.d345dcc8bf0012c19232a2d193a32928fc9c98ef.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/db1b1a44ce3cb041d5f412c3a023afc00daadd43/e2e/3976662b-7e18-4ef4-9397-f67d8c4c4bc5.md", "", "", "efbf9219-8bb8-4182-8079-84878fca9520.md")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/fbe7db58dccb6acbea8c051773a92bf3fe327b1d/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3976662b-7e18-4ef4-9397-f67d8c4c4bc5.7fed2e4b03b1d9d3d15c2b653214b755747046ba.zh-cn.xlf", "", "", "efbf9219-8bb8-4182-8079-84878fca9520.d345dcc8bf0012c19232a2d193a32928fc9c98ef.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/74cfad7bf15d4feed10d948ffefb164bb9c05f74/e2e/efbf9219-8bb8-4182-8079-84878fca9520.md", "", "", "3976662b-7e18-4ef4-9397-f67d8c4c4bc5.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/988e7166473fb56d35401be470c3d6e8703d49a9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/efbf9219-8bb8-4182-8079-84878fca9520.d345dcc8bf0012c19232a2d193a32928fc9c98ef.zh-cn.xlf", "", "", "3976662b-7e18-4ef4-9397-f67d8c4c4bc5.7fed2e4b03b1d9d3d15c2b653214b755747046ba.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/db1b1a44ce3cb041d5f412c3a023afc00daadd43/e2e/efbf9219-8bb8-4182-8079-84878fca9520.md", "", "", "3976662b-7e18-4ef4-9397-f67d8c4c4bc5.md")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/fbe7db58dccb6acbea8c051773a92bf3fe327b1d/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/efbf9219-8bb8-4182-8079-84878fca9520.d345dcc8bf0012c19232a2d193a32928fc9c98ef.zh-cn.xlf", "", "", "3976662b-7e18-4ef4-9397-f67d8c4c4bc5.7fed2e4b03b1d9d3d15c2b653214b755747046ba.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/74cfad7bf15d4feed10d948ffefb164bb9c05f74/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "efbf9219-8bb8-4182-8079-84878fca9520.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "efbf9219-8bb8-4182-8079-84878fca9520.d345dcc8bf0012c19232a2d193a32928fc9c98ef.de-de.xlf"
$ws.Range("D2").Value = "2016-03-08 20:53:55"
$ws.Range("E2").Value = "efbf9219-8bb8-4182-8079-84878fca9520.md"
$ws.Range("F2").Value = "efbf9219-8bb8-4182-8079-84878fca9520.d345dcc8bf0012c19232a2d193a32928fc9c98ef.de-de.xlf"
$ws.Range("G2").Value = "2016-03-08 20:54:53"
$ws.Range("H2").Value = "Include"

$ws.Range("A3").Value = "3976662b-7e18-4ef4-9397-f67d8c4c4bc5.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "3976662b-7e18-4ef4-9397-f67d8c4c4bc5.7fed2e4b03b1d9d3d15c2b653214b755747046ba.de-de.xlf"
$ws.Range("D3").Value = "2016-03-08 20:55:51"
$ws.Range("E3").Value = "3976662b-7e18-4ef4-9397-f67d8c4c4bc5.md"
$ws.Range("F3").Value = "3976662b-7e18-4ef4-9397-f67d8c4c4bc5.7fed2e4b03b1d9d3d15c2b653214b755747046ba.de-de.xlf"
$ws.Range("G3").Value = "2016-03-08 20:54:53"
$ws.Range("H3").Value = "Include"

$ws.Range("A4").Value = ".localization-config"
$ws.Range("B4").Value = "Not to be localized"
$ws.Range("D4").Value = "0001-01-01 00:00:00"
$ws.Range("G4").Value = "0001-01-01 00:00:00"
$ws.Range("H4").Value = "Ignored"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/74cfad7bf15d4feed10d948ffefb164bb9c05f74/e2e/3976662b-7e18-4ef4-9397-f67d8c4c4bc5.md", "", "", "efbf9219-8bb8-4182-8079-84878fca9520.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1d002a2f8369e6f6c50cd5cd9e2ea778f3cfae21/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3976662b-7e18-4ef4-9397-f67d8c4c4bc5.7fed2e4b03b1d9d3d15c2b653214b755747046ba.de-de.xlf", "", "", "efbf9219-8bb8-4182-8079-84878fca9520.d345dcc8bf0012c19232a2d193a32928fc9c98ef.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/7506e9fbf3e9066d6f5ee1cabf5260cde60187fe/e2e/3976662b-7e18-4ef4-9397-f67d8c4c4bc5.md", "", "", "efbf9219-8bb8-4182-8079-84878fca9520.md")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c2960a358dbe1035fc1015cc946dbf7ba16f29ce/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3976662b-7e18-4ef4-9397-f67d8c4c4bc5.7fed2e4b03b1d9d3d15c2b653214b755747046ba.de-de.xlf", "", "", "efbf9219-8bb8-4182-8079-84878fca9520.d345dcc8bf0012c19232a2d193a32928fc9c98ef.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/74cfad7bf15d4feed10d948ffefb164bb9c05f74/e2e/efbf9219-8bb8-4182-8079-84878fca9520.md", "", "", "3976662b-7e18-4ef4-9397-f67d8c4c4bc5.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1d002a2f8369e6f6c50cd5cd9e2ea778f3cfae21/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/efbf9219-8bb8-4182-8079-84878fca9520.d345dcc8bf0012c19232a2d193a32928fc9c98ef.de-de.xlf", "", "", "3976662b-7e18-4ef4-9397-f67d8c4c4bc5.7fed2e4b03b1d9d3d15c2b653214b755747046ba.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/7506e9fbf3e9066d6f5ee1cabf5260cde60187fe/e2e/efbf9219-8bb8-4182-8079-84878fca9520.md", "", "", "3976662b-7e18-4ef4-9397-f67d8c4c4bc5.md")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c2960a358dbe1035fc1015cc946dbf7ba16f29ce/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/efbf9219-8bb8-4182-8079-84878fca9520.d345dcc8bf0012c19232a2d193a32928fc9c98ef.de-de.xlf", "", "", "3976662b-7e18-4ef4-9397-f67d8c4c4bc5.7fed2e4b03b1d9d3d15c2b653214b755747046ba.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/74cfad7bf15d4feed10d948ffefb164bb9c05f74/.localization-config", "", "", ".localization-config")

$ws = $wb.Worksheets.Item("Overview")
$ws.Activate()
